$wb = $excel.ActiveWorkbook

# --- TestCases sheet: add new test case row ---
$wsTC = $wb.Worksheets.Item("TestCases")
$wsTC.Range("A7").Value = "CreateBillingTicketTest"
$wsTC.Range("B7").Value = "Y"
$wsTC.Range("B7").Select()

# --- Data sheet: add new test section with data ---
$wsData = $wb.Worksheets.Item("Data")

$wsData.Range("A26").Value = "CreateBillingTicketTest"
$wsData.Range("A26").Style = "Normal"

$wsData.Range("A27").Value = "Runmode"
$wsData.Range("B27").Value = "Result"
$wsData.Range("C27").Value = "Browser"

$wsData.Range("D27").Value = "Account # For Billing Ticket"
$wsData.Range("D28").Value = "'200069083"

$wsData.Range("E27").Value = "Billing ticket title"
$wsData.Range("E28").Value = "testtitle"

$wsData.Range("F27").Value = "What is issue"
$wsData.Range("F28").Value = "bill not generating"

$wsData.Range("G27").Value = "Col4"

$wsData.Range("A28").Value = "Y"
$wsData.Range("C28").Value = "Mozilla"

$wsData.Select()
$wsData.Range("D30").Select()
